$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 404, pushing existing rows 404:523 down to 405:524
$ws.Rows(404).Insert()

# Copy the content of the now-shifted original row (now at 405) into the new row 404
$ws.Range("A405:R405").Copy() | Out-Null
$ws.Range("A404").PasteSpecial() | Out-Null

# Update the two changed values for the newly inserted row
$ws.Range("D404").Value = 45093
$ws.Range("J404").Value = 55
